$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for columns I and J (matching style of existing header cells)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the I0 / IF data columns for rows 2-64
$iValues = @(8,9,4,8,9,7,8,9,9,6,7,7,12,8,5,6,7,9,8,6,4,8,7,6,5,7,8,6,7,8,7,6,6,8,7,7,6,6,8,9,8,9,8,8,6,6,7,7,7,7,6,6,6,6,5,6,5,4,9,6,7,4,7)
$jValues = @(8,9,4,8,9,7,8,9,9,7,7,7,12,8,6,6,8,9,9,7,6,8,7,6,7,7,8,6,8,9,8,7,8,8,9,8,7,6,8,9,9,10,8,8,7,7,8,7,7,8,7,7,6,7,6,6,5,5,9,6,7,4,7)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
